$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column D, shifting existing D:H to E:I
$ws.Columns("D:D").Insert()

# Set header for new column D
$ws.Range("D1").Value = "id"

# Set "NL" value for each data row in the new column D
$ws.Range("D2:D11").Value = "NL"

# Match the final selection shown in the saved file
$ws.Range("D11").Select() | Out-Null
